$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new "Movimiento" column at F, shifting existing columns F..Y one
# position to the right (-> G..Z), without touching the existing column
# width/style definitions. We do this by copying each column range (values +
# formats) from the rightmost column down to F into the column immediately
# to its right, working right-to-left so a destination is always written
# before it is used as a source on the next step.

$pairs = @(
    @{src = "Y"; dst = "Z"},
    @{src = "X"; dst = "Y"},
    @{src = "W"; dst = "X"},
    @{src = "V"; dst = "W"},
    @{src = "U"; dst = "V"},
    @{src = "T"; dst = "U"},
    @{src = "S"; dst = "T"},
    @{src = "R"; dst = "S"},
    @{src = "Q"; dst = "R"},
    @{src = "P"; dst = "Q"},
    @{src = "O"; dst = "P"},
    @{src = "N"; dst = "O"},
    @{src = "M"; dst = "N"},
    @{src = "L"; dst = "M"},
    @{src = "K"; dst = "L"},
    @{src = "J"; dst = "K"},
    @{src = "I"; dst = "J"},
    @{src = "H"; dst = "I"},
    @{src = "G"; dst = "H"},
    @{src = "F"; dst = "G"}
)

foreach ($pair in $pairs) {
    $srcAddr = [string]$pair.src + "1:" + [string]$pair.src + "11"
    $dstAddr = [string]$pair.dst + "1"
    $ws.Range($srcAddr).Copy($ws.Range($dstAddr))
    $ws.Range($srcAddr).ClearContents()
}

# Column F is now free for the new "Movimiento" header; the data rows below
# it stay blank, matching the target layout.
$ws.Range("F1").Value = "Movimiento"

$wb.Save()
